$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current "Buyer Table" block (row 9) to make
# room for the Farmer Table's new "Username" / "Password" login fields.
# This shifts the Buyer Table and Admin Table blocks down by two rows each,
# including merged cells and the blank separator rows, just like Excel's own
# Rows.Insert.
$ws.Rows("9:10").Insert()

# Give the two new rows the same bordered look as the rest of the table by
# copying the format from the row right above (Address / row 8).
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "Username" and "Password" property rows.
$ws.Range("A9").Value = "Username"
$ws.Range("B9").Value = "Varchar"
$ws.Range("C9").Value = "Unique Key"
$ws.Range("D9").Value = "Not Null"

$ws.Range("A10").Value = "Password"
$ws.Range("B10").Value = "Varchar"
$ws.Range("D10").Value = "Not Null"

# Those two rows are not fully complete (C10 stays blank, same as the other
# property rows), matching the surrounding table's look.
$ws.Range("A9").Interior.Color = 16777215
$ws.Range("A9").Interior.Pattern = 1
$ws.Range("B9").Interior.Color = 16777215
$ws.Range("B9").Interior.Pattern = 1
$ws.Range("D9").Interior.Color = 16777215
$ws.Range("D9").Interior.Pattern = 1
$ws.Range("A10").Interior.Color = 16777215
$ws.Range("A10").Interior.Pattern = 1
$ws.Range("B10").Interior.Color = 16777215
$ws.Range("B10").Interior.Pattern = 1

# Restore the view: the author had scrolled down and was working around
# C10 when the file was saved.
$ws.Range("C10").Select()
